$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Range("H5").Value = "Manage Class"
$ws.Range("H6").Value = "Dashboard"

$ws.Activate()
$ws.Range("H2").Select()
